$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 17.04324350258893
$ws.Cells.Item(2, 3).Value = 5.598243218831244
$ws.Cells.Item(2, 4).Value = 11.64155665218456
$ws.Cells.Item(2, 5).Value = 11.26414711091073
$ws.Cells.Item(2, 6).Value = 57.93622936338958
$ws.Cells.Item(2, 11).Value = 13.41342992578843
$ws.Cells.Item(2, 12).Value = 10.20684388430188
$ws.Cells.Item(2, 13).Value = 16.55792841314695
$ws.Cells.Item(3, 2).Value = 17.02031697725667
$ws.Cells.Item(3, 3).Value = 5.506379353198493
$ws.Cells.Item(3, 4).Value = 11.50265263352384
$ws.Cells.Item(3, 5).Value = 11.25634870994377
$ws.Cells.Item(3, 6).Value = 56.83751660945563
$ws.Cells.Item(3, 11).Value = 13.4276402657267
$ws.Cells.Item(3, 12).Value = 10.21989284324145
$ws.Cells.Item(3, 13).Value = 16.59012411845315
$ws.Cells.Item(4, 2).Value = 17.0124793732986
$ws.Cells.Item(4, 3).Value = 5.447486261030732
$ws.Cells.Item(4, 4).Value = 11.41534165416848
$ws.Cells.Item(4, 5).Value = 11.25220294936749
$ws.Cells.Item(4, 6).Value = 56.15357939601261
$ws.Cells.Item(4, 11).Value = 13.44154486887179
$ws.Cells.Item(4, 12).Value = 10.22934092039317
$ws.Cells.Item(4, 13).Value = 16.6134901488376
$ws.Cells.Item(5, 2).Value = 17.01085662783823
$ws.Cells.Item(5, 3).Value = 5.422862620463366
$ws.Cells.Item(5, 4).Value = 11.37926595019279
$ws.Cells.Item(5, 5).Value = 11.25067545713118
$ws.Cells.Item(5, 6).Value = 55.87278340861413
$ws.Cells.Item(5, 11).Value = 13.44850991437816
$ws.Cells.Item(5, 12).Value = 10.23355223626543
$ws.Cells.Item(5, 13).Value = 16.62391557373416
$ws.Cells.Item(6, 2).Value = 17.01068209526821
$ws.Cells.Item(6, 3).Value = 5.418736228742923
$ws.Cells.Item(6, 4).Value = 11.37324604819275
$ws.Cells.Item(6, 5).Value = 11.25043160931813
$ws.Cells.Item(6, 6).Value = 55.82603915182034
$ws.Cells.Item(6, 11).Value = 13.44974477585311
$ws.Cells.Item(6, 12).Value = 10.23427333545691
$ws.Cells.Item(6, 13).Value = 16.62570124558386
$ws.Cells.Item(7, 2).Value = 17.01245112535912
$ws.Cells.Item(7, 3).Value = 5.447156703327102
$ws.Cells.Item(7, 4).Value = 11.41485711454436
$ws.Cells.Item(7, 5).Value = 11.25218169289311
$ws.Cells.Item(7, 6).Value = 56.14980058278244
$ws.Cells.Item(7, 11).Value = 13.44163354862211
$ws.Cells.Item(7, 12).Value = 10.22939625335341
$ws.Cells.Item(7, 13).Value = 16.61362709284701
$ws.Cells.Item(8, 2).Value = 17.03404530072966
$ws.Cells.Item(8, 3).Value = 5.567085822275061
$ws.Cells.Item(8, 4).Value = 11.59409021939295
$ws.Cells.Item(8, 5).Value = 11.26132510850146
$ws.Cells.Item(8, 6).Value = 57.55948079948023
$ws.Cells.Item(8, 11).Value = 13.41725279354989
$ws.Cells.Item(8, 12).Value = 10.21104518034349
$ws.Cells.Item(8, 13).Value = 16.56828217153397
$ws.Cells.Item(9, 2).Value = 17.12572081233726
$ws.Cells.Item(9, 3).Value = 5.782472406115656
$ws.Cells.Item(9, 4).Value = 11.92895759623391
$ws.Cells.Item(9, 5).Value = 11.28433671937169
$ws.Cells.Item(9, 6).Value = 60.2389679948549
$ws.Cells.Item(9, 11).Value = 13.41065471049142
$ws.Cells.Item(9, 12).Value = 10.18644900253945
$ws.Cells.Item(9, 13).Value = 16.50795044169884
$ws.Cells.Item(10, 2).Value = 17.22279296910945
$ws.Cells.Item(10, 3).Value = 5.928603314562781
$ws.Cells.Item(10, 4).Value = 12.16414639955432
$ws.Cells.Item(10, 5).Value = 11.30431936471294
$ws.Cells.Item(10, 6).Value = 62.14191684388643
$ws.Cells.Item(10, 11).Value = 13.43102447109744
$ws.Cells.Item(10, 12).Value = 10.17531611461027
$ws.Cells.Item(10, 13).Value = 16.48110324242463
$ws.Cells.Item(11, 2).Value = 17.27329414357659
$ws.Cells.Item(11, 3).Value = 5.99244604599351
$ws.Cells.Item(11, 4).Value = 12.26864198823409
$ws.Cells.Item(11, 5).Value = 11.31407264080164
$ws.Cells.Item(11, 6).Value = 62.99062873151685
$ws.Cells.Item(11, 11).Value = 13.44576468010347
$ws.Cells.Item(11, 12).Value = 10.17175581775974
$ws.Cells.Item(11, 13).Value = 16.4726915019759
$ws.Cells.Item(12, 2).Value = 17.29331749727385
$ws.Cells.Item(12, 3).Value = 6.016242570599236
$ws.Cells.Item(12, 4).Value = 12.30784141913639
$ws.Cells.Item(12, 5).Value = 11.31786073399588
$ws.Cells.Item(12, 6).Value = 63.30936420109315
$ws.Cells.Item(12, 11).Value = 13.45213155604446
$ws.Cells.Item(12, 12).Value = 10.17062362114826
$ws.Cells.Item(12, 13).Value = 16.47005291274659
$ws.Cells.Item(13, 2).Value = 17.28896531240389
$ws.Cells.Item(13, 3).Value = 6.011134443665326
$ws.Cells.Item(13, 4).Value = 12.29941580214344
$ws.Cells.Item(13, 5).Value = 11.31704069693429
$ws.Cells.Item(13, 6).Value = 63.24084012529908
$ws.Cells.Item(13, 11).Value = 13.45072546345457
$ws.Cells.Item(13, 12).Value = 10.17085785784898
$ws.Cells.Item(13, 13).Value = 16.47059686324128
$ws.Cells.Item(14, 2).Value = 17.27492352342292
$ws.Cells.Item(14, 3).Value = 5.994411406238529
$ws.Cells.Item(14, 4).Value = 12.27187444890734
$ws.Cells.Item(14, 5).Value = 11.31438239442565
$ws.Cells.Item(14, 6).Value = 63.01690565868736
$ws.Cells.Item(14, 11).Value = 13.44627276844013
$ws.Cells.Item(14, 12).Value = 10.17165834388433
$ws.Cells.Item(14, 13).Value = 16.47246346587112
$ws.Cells.Item(15, 2).Value = 17.26643928752914
$ws.Cells.Item(15, 3).Value = 5.984118644509675
$ws.Cells.Item(15, 4).Value = 12.25495593059482
$ws.Cells.Item(15, 5).Value = 11.31276642645524
$ws.Cells.Item(15, 6).Value = 62.87938759925667
$ws.Cells.Item(15, 11).Value = 13.44364752440363
$ws.Cells.Item(15, 12).Value = 10.17217678642309
$ws.Cells.Item(15, 13).Value = 16.47367801780909
$ws.Cells.Item(16, 2).Value = 17.21961920555282
$ws.Cells.Item(16, 3).Value = 5.924377994781813
$ws.Cells.Item(16, 4).Value = 12.15726614155798
$ws.Cells.Item(16, 5).Value = 11.30369525667162
$ws.Cells.Item(16, 6).Value = 62.08609170894204
$ws.Cells.Item(16, 11).Value = 13.43017116185261
$ws.Cells.Item(16, 12).Value = 10.17557903236925
$ws.Cells.Item(16, 13).Value = 16.48172947707368
$ws.Cells.Item(17, 2).Value = 17.19251279670273
$ws.Cells.Item(17, 3).Value = 5.88705411148388
$ws.Cells.Item(17, 4).Value = 12.09668923011385
$ws.Cells.Item(17, 5).Value = 11.29829979518066
$ws.Cells.Item(17, 6).Value = 61.59493157174713
$ws.Cells.Item(17, 11).Value = 13.42330473408409
$ws.Cells.Item(17, 12).Value = 10.1780512580133
$ws.Cells.Item(17, 13).Value = 16.48764254134748
$ws.Cells.Item(18, 2).Value = 17.17751975015463
$ws.Cells.Item(18, 3).Value = 5.865338606280932
$ws.Cells.Item(18, 4).Value = 12.06161305155901
$ws.Cells.Item(18, 5).Value = 11.29525891880168
$ws.Cells.Item(18, 6).Value = 61.31084604886938
$ws.Cells.Item(18, 11).Value = 13.41987077039057
$ws.Cells.Item(18, 12).Value = 10.17961481094856
$ws.Cells.Item(18, 13).Value = 16.49140137106491
$ws.Cells.Item(19, 2).Value = 17.17254639251622
$ws.Cells.Item(19, 3).Value = 5.857943575223054
$ws.Cells.Item(19, 4).Value = 12.04969706201364
$ws.Cells.Item(19, 5).Value = 11.29424007413741
$ws.Cells.Item(19, 6).Value = 61.21439444614592
$ws.Cells.Item(19, 11).Value = 13.4187966583662
$ws.Cells.Item(19, 12).Value = 10.18016852960086
$ws.Cells.Item(19, 13).Value = 16.49273548884604
$ws.Cells.Item(20, 2).Value = 17.19533653039928
$ws.Cells.Item(20, 3).Value = 5.891052954026501
$ws.Cells.Item(20, 4).Value = 12.10316204801333
$ws.Cells.Item(20, 5).Value = 11.29886769215234
$ws.Cells.Item(20, 6).Value = 61.64738186772291
$ws.Cells.Item(20, 11).Value = 13.42398234417835
$ws.Cells.Item(20, 12).Value = 10.17777343278352
$ws.Cells.Item(20, 13).Value = 16.48697605377451
$ws.Cells.Item(21, 2).Value = 17.27902362787177
$ws.Cells.Item(21, 3).Value = 5.999333667492571
$ws.Cells.Item(21, 4).Value = 12.27997417072228
$ws.Cells.Item(21, 5).Value = 11.31516063617305
$ws.Cells.Item(21, 6).Value = 63.08275430814363
$ws.Cells.Item(21, 11).Value = 13.44755934803398
$ws.Cells.Item(21, 12).Value = 10.17141736176456
$ws.Cells.Item(21, 13).Value = 16.47190036066096
$ws.Cells.Item(22, 2).Value = 17.33895495605253
$ws.Cells.Item(22, 3).Value = 6.06789012901113
$ws.Cells.Item(22, 4).Value = 12.39336452217582
$ws.Cells.Item(22, 5).Value = 11.32636091869175
$ws.Cells.Item(22, 6).Value = 64.00529653413314
$ws.Cells.Item(22, 11).Value = 13.46754258186453
$ws.Cells.Item(22, 12).Value = 10.16852220789103
$ws.Cells.Item(22, 13).Value = 16.46523432773504
$ws.Cells.Item(23, 2).Value = 17.30649381559633
$ws.Cells.Item(23, 3).Value = 6.031502774656835
$ws.Cells.Item(23, 4).Value = 12.33304815165963
$ws.Cells.Item(23, 5).Value = 11.3203328297572
$ws.Cells.Item(23, 6).Value = 63.51440988542011
$ws.Cells.Item(23, 11).Value = 13.45645957347591
$ws.Cells.Item(23, 12).Value = 10.1699523224801
$ws.Cells.Item(23, 13).Value = 16.46850053422873
$ws.Cells.Item(24, 2).Value = 17.19405807921818
$ws.Cells.Item(24, 3).Value = 5.88924587940117
$ws.Cells.Item(24, 4).Value = 12.10023646396095
$ws.Cells.Item(24, 5).Value = 11.29861075575499
$ws.Cells.Item(24, 6).Value = 61.62367440177439
$ws.Cells.Item(24, 11).Value = 13.42367439656815
$ws.Cells.Item(24, 12).Value = 10.17789859450032
$ws.Cells.Item(24, 13).Value = 16.48727625353628
$ws.Cells.Item(25, 2).Value = 17.0956686428771
$ws.Cells.Item(25, 3).Value = 5.726329247762829
$ws.Cells.Item(25, 4).Value = 11.84021997836694
$ws.Cells.Item(25, 5).Value = 11.27756955653279
$ws.Cells.Item(25, 6).Value = 59.52466198294128
$ws.Cells.Item(25, 11).Value = 13.40801379871049
$ws.Cells.Item(25, 12).Value = 10.19188392218093
$ws.Cells.Item(25, 13).Value = 16.52120511107039
